$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Peaches" samples (and a couple of other blocks) were originally
# logged with an "Unknown" collection location; they have since been
# identified as coming from Colorado. There is a single shared-string
# entry for the literal text "Unknown" in this workbook, so a sheet-wide
# replace updates every cell that pointed at it (columns F457:F552,
# F1431:F1526 and F2405:F2500) without touching any other text.
$ws.Cells.Replace("Unknown", "Colorado")

# Reflect the author's final on-screen selection/scroll state: the newly
# edited location column (F2405:F2500) ends up selected, with F2405 as
# the active cell.
$ws.Range("F2405:F2500").Select()
